$wb = $excel.ActiveWorkbook

# Sheet "Overview": Status columns for zh-cn (B) and de-de (C)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"
$wsOverview.Range("B5").Value = "In Translation"
$wsOverview.Range("C5").Value = "In Translation"

# Sheet "zh-cn": Status column (B)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B4").Value = "In Translation"
$wsZhCn.Range("B5").Value = "In Translation"

# Sheet "de-de": Status column (B)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B4").Value = "In Translation"
$wsDeDe.Range("B5").Value = "In Translation"
